# Move the first 40 "phone" rows (A:B) into the processed-log columns (D:E),
# shifted down by one row so D1 can hold a marker header, and shift the
# remaining (unprocessed) rows up so column A:B restarts at row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$processedCount = 40

# Read every existing A/B pair into arrays before we start overwriting cells.
$aVals = @()
$bVals = @()
for ($i = 1; $i -le $lastRow; $i++) {
    $aVals += ,$ws.Cells.Item($i, 1).Value()
    $bVals += ,$ws.Cells.Item($i, 2).Value()
}

# Marker header for the processed-log column.
$ws.Cells.Item(1, 4).Value = "PROCESSED ALREADY PHONES IN THIS COLUMN"

# Old rows 1..40 (A,B) -> new D2:E41 (the "already processed" log).
for ($i = 1; $i -le $processedCount; $i++) {
    $ws.Cells.Item($i + 1, 4).Value = $aVals[$i - 1]
    $ws.Cells.Item($i + 1, 5).Value = $bVals[$i - 1]
}

# Old rows 41..lastRow (A,B) -> new A1:B(lastRow-40) (still-unprocessed phones).
$newLastRow = $lastRow - $processedCount
for ($i = $processedCount + 1; $i -le $lastRow; $i++) {
    $newRow = $i - $processedCount
    $ws.Cells.Item($newRow, 1).Value = $aVals[$i - 1]
    $ws.Cells.Item($newRow, 2).Value = $bVals[$i - 1]
}

# Clear the now-stale tail of columns A:B (rows that used to hold data past
# the new shortened range).
for ($i = $newLastRow + 1; $i -le $lastRow; $i++) {
    $ws.Cells.Item($i, 1).Value = ""
    $ws.Cells.Item($i, 2).Value = ""
}

$ws.Range("A2").Select() | Out-Null
